# Apply the changes described by the diff:
#  1. Insert a new sheet "AddGeneIdSource" right after "RepositoryIdSource"
#     with columns: sequence_file_url, addgene_sequence_type, repository_name,
#     repository_id, input, output, type, id
#     plus two list data validations (columns B and C).
#  2. Insert a new sheet "CRISPRSource" right after "RestrictionAndLigationSource"
#     with columns: guides, circular, assembly, input, output, type, id

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. AddGeneIdSource
# ---------------------------------------------------------------------------
$afterSheet1 = $wb.Worksheets.Item("RepositoryIdSource")
$addGeneIdSource = $wb.Worksheets.Add($null, $afterSheet1)
$addGeneIdSource.Name = "AddGeneIdSource"

$addGeneIdHeaders = @(
    "sequence_file_url",
    "addgene_sequence_type",
    "repository_name",
    "repository_id",
    "input",
    "output",
    "type",
    "id"
)

for ($i = 0; $i -lt $addGeneIdHeaders.Length; $i++) {
    $addGeneIdSource.Cells.Item(1, $i + 1).Value = $addGeneIdHeaders[$i]
}

# Data validation on column B (addgene_sequence_type)
$colB = $addGeneIdSource.Range("B2:B1048576")
$colB.Validation.Add(3, 1, $null, '"depositor-full,addgene-full"')
$colB.Validation.ShowInput = $false
$colB.Validation.ShowError = $false

# Data validation on column C (repository_name)
$colC = $addGeneIdSource.Range("C2:C1048576")
$colC.Validation.Add(3, 1, $null, '"addgene,genbank"')
$colC.Validation.ShowInput = $false
$colC.Validation.ShowError = $false

# ---------------------------------------------------------------------------
# 2. CRISPRSource
# ---------------------------------------------------------------------------
$afterSheet2 = $wb.Worksheets.Item("RestrictionAndLigationSource")
$crisprSource = $wb.Worksheets.Add($null, $afterSheet2)
$crisprSource.Name = "CRISPRSource"

$crisprHeaders = @(
    "guides",
    "circular",
    "assembly",
    "input",
    "output",
    "type",
    "id"
)

for ($i = 0; $i -lt $crisprHeaders.Length; $i++) {
    $crisprSource.Cells.Item(1, $i + 1).Value = $crisprHeaders[$i]
}
